$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.400.02"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "3.642.35"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'196.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.67%  "
$ws.Range("D6").Value = "'582.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "3.638.76"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "'0.620"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  +6.63%  "
$ws.Range("D12").Value = "'56.64"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.31%  "
$ws.Range("D13").Value = "'0.0000293"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +14.91%  "
$ws.Range("D14").Value = "'10.14"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "4.233.35"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "3.651.85"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "'12.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("D19").Value = "68.439.95"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "'18.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "'403.38"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").Value = "'13.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +26.44%  "
$ws.Range("D24").Value = "'4.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'86.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "'2.96"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.32%  "
$ws.Range("D27").Value = "'12.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("D28").Value = "'3.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.67%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "'8.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +20.29%  "
$ws.Range("D31").Value = "'9.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").Value = "'31.81"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").Value = "'706.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +19.10%  "
$ws.Range("D34").Value = "'12.25"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  +5.77%  "
$ws.Range("D36").Value = "'64.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'42.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").Value = "'0.426"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +13.54%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "0.0₃0790"
$ws.Range("E40").Value = "  +5.88%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.88"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +19.24%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.138"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.233.80"
$ws.Range("E43").Value = "  +20.45%  "
$ws.Range("D44").Value = "'3.14"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +13.16%  "
$ws.Range("D45").Value = "'3.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +35.74%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'0.0423"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "'8.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.24%  "
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").Value = "'3.10"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'2.63"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.78%  "
